# Update the "expected experiment duration" sentence on the consent slide.
#
# Original (3 runs in the paragraph):
#   run1: "...משך הניסוי הצפוי הוא כ"
#   run2: "40"
#   run3: " דקות."
#
# New (3 runs in the paragraph):
#   run1: "...משך הניסוי הצפוי "
#   run2: "הוא כ- 15 דקות"
#   run3: "."

$p = $ppt.ActivePresentation

# Find the slide/shape that holds the consent paragraph mentioning "40" minutes.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $t = $shape.TextFrame.TextRange.Text
            if ($t -like "*הצפוי הוא כ40 דקות*") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -eq $null) {
    # Fallback to the known location (slide 3, "Content Placeholder 2").
    $targetSlide = $p.Slides.Item(3)
    $targetShape = $targetSlide.Shapes.Item("Content Placeholder 2")
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text

# Original run texts (used only to locate the exact character offsets of
# each run in the current document).
$origRun1Text = "ניתן להפסיק בכל שלב את השתתפותך בניסוי. אך אנו נוכל להשתמש רק במידע מניסויים שהושלמו לצורך המחקר. משך הניסוי הצפוי הוא כ"
$origRun2Text = "40"
$origRun3Text = " דקות."

$run1Start0 = $fullText.IndexOf($origRun1Text + $origRun2Text + $origRun3Text)   # 0-based

if ($run1Start0 -ge 0) {
    $run1Start = $run1Start0 + 1                                                 # 1-based
    $run1Len = $origRun1Text.Length
    $run2Start = $run1Start + $run1Len
    $run2Len = $origRun2Text.Length
    $run3Start = $run2Start + $run2Len
    $run3Len = $origRun3Text.Length

    # New text for each run.
    $newRun1Text = "ניתן להפסיק בכל שלב את השתתפותך בניסוי. אך אנו נוכל להשתמש רק במידע מניסויים שהושלמו לצורך המחקר. משך הניסוי הצפוי "
    $newRun2Text = "הוא כ- 15 דקות"
    $newRun3Text = "."

    # Apply the edits back-to-front so earlier (lower) offsets stay valid as
    # the text lengths change.
    $r3 = $tr.Characters($run3Start, $run3Len)
    $r3.Text = $newRun3Text

    $r2 = $tr.Characters($run2Start, $run2Len)
    $r2.Text = $newRun2Text

    $r1 = $tr.Characters($run1Start, $run1Len)
    $r1.Text = $newRun1Text
}
